# Fruta / hortaliza, semanal
# Insert a new weekly record row for "Camote" at row 120 of the Zapallo
# sheet (Terminal Hortofrutícola Agro Chillán), pushing the existing
# rows 120-138 down to 121-139.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 120; this shifts rows 120:138
# down to 121:139 and carries the row formatting down with them.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new weekly observation.
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44637
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112045
$ws.Cells.Item(120, 7).Value = "Zapallo"
$ws.Cells.Item(120, 8).Value = "Camote"
$ws.Cells.Item(120, 9).Value = "1a (cosecha)"
$ws.Cells.Item(120, 10).Value = 400
$ws.Cells.Item(120, 11).Value = 300
$ws.Cells.Item(120, 12).Value = 350
$ws.Cells.Item(120, 13).Value = 325
$ws.Cells.Item(120, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(120, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(120, 16).Value = 325
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"
